$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: Status column C2 + Latest Handoff Datetime H2
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-23 12:38:57"

# de-de sheet: Status column C2 + Latest Handoff Datetime H2
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-23 12:39:07"

# Overview sheet: zh-cn status E2, de-de status F2, Latest HO Xliff Generate Date G2
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-23 12:39:07"

# Column widths autofit to match new (longer) text
$overview.Columns.Item(5).AutoFit() | Out-Null
$overview.Columns.Item(6).AutoFit() | Out-Null
$zhcn.Columns.Item(3).AutoFit() | Out-Null
$dede.Columns.Item(3).AutoFit() | Out-Null
